# GitHub Actions "Updated symbol list" refresh (Tue Jan 10 07:09:53 UTC 2023).
#
# Sheet1 is a scraped coin-ranking table. Columns D (Price), E (Volume(1h))
# and G (Hora) are stored as literal text (the sheet uses inlineStr cells,
# not numbers), so every write below is done as text: a leading apostrophe
# stops Excel from re-typing a numeric-looking string ("274.66", "-1.19%",
# "7", ...) as a Double/percentage, and the cell's Style is captured and
# restored around the write so the implicit "number stored as text" quote
# -prefix style never gets attached to the cell (keeping formatting as it
# was before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: new Price / Volume(1h) / Hora text for that cell. Rows whose
# Price+Volume are still the placeholder "--" / "--%" only get Hora bumped.
$updates = @(
    @{ Cell = 'D2'; Text = '274.66' },
    @{ Cell = 'E2'; Text = '-1.19%' },
    @{ Cell = 'G2'; Text = '7' },
    @{ Cell = 'D3'; Text = '26.81' },
    @{ Cell = 'E3'; Text = '-2.19%' },
    @{ Cell = 'G3'; Text = '7' },
    @{ Cell = 'D4'; Text = '4.859' },
    @{ Cell = 'E4'; Text = '1.05%' },
    @{ Cell = 'G4'; Text = '7' },
    @{ Cell = 'D5'; Text = '0.06327' },
    @{ Cell = 'E5'; Text = '1.32%' },
    @{ Cell = 'G5'; Text = '7' },
    @{ Cell = 'D6'; Text = '6.873' },
    @{ Cell = 'E6'; Text = '-0.46%' },
    @{ Cell = 'G6'; Text = '7' },
    @{ Cell = 'D7'; Text = '3.324' },
    @{ Cell = 'E7'; Text = '1.38%' },
    @{ Cell = 'G7'; Text = '7' },
    @{ Cell = 'D8'; Text = '1.262' },
    @{ Cell = 'E8'; Text = '33.87%' },
    @{ Cell = 'G8'; Text = '7' },
    @{ Cell = 'D9'; Text = '0.8680' },
    @{ Cell = 'E9'; Text = '-1.44%' },
    @{ Cell = 'G9'; Text = '7' },
    @{ Cell = 'D10'; Text = '0.1517' },
    @{ Cell = 'E10'; Text = '4.64%' },
    @{ Cell = 'G10'; Text = '7' },
    @{ Cell = 'D11'; Text = '0.04997' },
    @{ Cell = 'E11'; Text = '-4.54%' },
    @{ Cell = 'G11'; Text = '7' },
    @{ Cell = 'D12'; Text = '0.07453' },
    @{ Cell = 'E12'; Text = '2.30%' },
    @{ Cell = 'G12'; Text = '7' },
    @{ Cell = 'D13'; Text = '0.02952' },
    @{ Cell = 'E13'; Text = '-6.55%' },
    @{ Cell = 'G13'; Text = '7' },
    @{ Cell = 'D14'; Text = '0.08990' },
    @{ Cell = 'E14'; Text = '-0.78%' },
    @{ Cell = 'G14'; Text = '7' },
    @{ Cell = 'D15'; Text = '0.001571' },
    @{ Cell = 'E15'; Text = '0.54%' },
    @{ Cell = 'G15'; Text = '7' },
    @{ Cell = 'D16'; Text = '0.0006297' },
    @{ Cell = 'E16'; Text = '0.24%' },
    @{ Cell = 'G16'; Text = '7' },
    @{ Cell = 'D17'; Text = '0.005977' },
    @{ Cell = 'E17'; Text = '4.07%' },
    @{ Cell = 'G17'; Text = '7' },
    @{ Cell = 'D18'; Text = '3.448' },
    @{ Cell = 'E18'; Text = '-0.16%' },
    @{ Cell = 'G18'; Text = '7' },
    @{ Cell = 'D19'; Text = '2.272' },
    @{ Cell = 'E19'; Text = '-0.54%' },
    @{ Cell = 'G19'; Text = '7' },
    @{ Cell = 'D20'; Text = '0.3126' },
    @{ Cell = 'E20'; Text = '1.03%' },
    @{ Cell = 'G20'; Text = '7' },
    @{ Cell = 'D21'; Text = '0.1330' },
    @{ Cell = 'E21'; Text = '2.81%' },
    @{ Cell = 'G21'; Text = '7' },
    @{ Cell = 'D22'; Text = '3.916' },
    @{ Cell = 'E22'; Text = '1.22%' },
    @{ Cell = 'G22'; Text = '7' },
    @{ Cell = 'D23'; Text = '0.04365' },
    @{ Cell = 'E23'; Text = '0.83%' },
    @{ Cell = 'G23'; Text = '7' },
    @{ Cell = 'D24'; Text = '0.001175' },
    @{ Cell = 'E24'; Text = '-0.10%' },
    @{ Cell = 'G24'; Text = '7' },
    @{ Cell = 'D25'; Text = '0.004246' },
    @{ Cell = 'E25'; Text = '-0.51%' },
    @{ Cell = 'G25'; Text = '7' },
    @{ Cell = 'E26'; Text = '-0.09%' },
    @{ Cell = 'G26'; Text = '7' },
    @{ Cell = 'E27'; Text = '-0.34%' },
    @{ Cell = 'G27'; Text = '7' },
    @{ Cell = 'G28'; Text = '7' },
    @{ Cell = 'G29'; Text = '7' },
    @{ Cell = 'G30'; Text = '7' },
    @{ Cell = 'G31'; Text = '7' },
    @{ Cell = 'G32'; Text = '7' },
    @{ Cell = 'G33'; Text = '7' },
    @{ Cell = 'G34'; Text = '7' },
    @{ Cell = 'G35'; Text = '7' },
    @{ Cell = 'G36'; Text = '7' },
    @{ Cell = 'G37'; Text = '7' },
    @{ Cell = 'G38'; Text = '7' },
    @{ Cell = 'G39'; Text = '7' },
    @{ Cell = 'D40'; Text = '0.04045' },
    @{ Cell = 'E40'; Text = '0.38%' },
    @{ Cell = 'G40'; Text = '7' },
    @{ Cell = 'D41'; Text = '0.006701' },
    @{ Cell = 'E41'; Text = '4.80%' },
    @{ Cell = 'G41'; Text = '7' },
    @{ Cell = 'D42'; Text = '0.1166' },
    @{ Cell = 'E42'; Text = '1.08%' },
    @{ Cell = 'G42'; Text = '7' },
    @{ Cell = 'D43'; Text = '0.002090' },
    @{ Cell = 'E43'; Text = '-2.08%' },
    @{ Cell = 'G43'; Text = '7' },
    @{ Cell = 'D44'; Text = '0.01068' },
    @{ Cell = 'E44'; Text = '-11.42%' },
    @{ Cell = 'G44'; Text = '7' },
    @{ Cell = 'D45'; Text = '0.00005286' },
    @{ Cell = 'E45'; Text = '3.88%' },
    @{ Cell = 'G45'; Text = '7' },
    @{ Cell = 'E46'; Text = '-33.04%' },
    @{ Cell = 'G46'; Text = '7' },
    @{ Cell = 'D47'; Text = '1.486' },
    @{ Cell = 'E47'; Text = '-37.45%' },
    @{ Cell = 'G47'; Text = '7' },
    @{ Cell = 'G48'; Text = '7' },
    @{ Cell = 'G49'; Text = '7' },
    @{ Cell = 'G50'; Text = '7' },
    @{ Cell = 'G51'; Text = '7' }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $savedStyle = $cell.Style
    $cell.Value = "'" + $update.Text
    $cell.Style = $savedStyle
}
